$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Note: "170÷9=18, 8" is both an old value (last replacement below) and a
# new value (first replacement below). Perform the replacement that
# consumes "170÷9=18, 8" as its search text FIRST, before the one that
# creates a new "170÷9=18, 8" string, so the Find doesn't accidentally
# match the freshly-written text.
Replace-Text "170÷9=18, 8" "146÷5=29, 1"

Replace-Text "563÷4=140, 3" "170÷9=18, 8"
Replace-Text "174÷7=24, 6" "575÷2=287, 1"
Replace-Text "507÷2=253, 1" "399÷6=66, 3"
Replace-Text "733÷9=81, 4" "413÷9=45, 8"
Replace-Text "555÷6=92, 3" "258÷5=51, 3"
Replace-Text "691÷6=115, 1" "154÷5=30, 4"
Replace-Text "852÷3=284, 0" "425÷6=70, 5"
Replace-Text "239÷5=47, 4" "505÷8=63, 1"
Replace-Text "781÷8=97, 5" "699÷6=116, 3"
Replace-Text "640÷9=71, 1" "137÷4=34, 1"
Replace-Text "875÷9=97, 2" "499÷3=166, 1"
Replace-Text "141÷9=15, 6" "685÷3=228, 1"
Replace-Text "318÷8=39, 6" "844÷3=281, 1"
Replace-Text "782÷3=260, 2" "581÷4=145, 1"
Replace-Text "976÷8=122, 0" "978÷4=244, 2"
Replace-Text "850÷5=170, 0" "785÷2=392, 1"
Replace-Text "459÷6=76, 3" "603÷3=201, 0"
Replace-Text "723÷7=103, 2" "270÷4=67, 2"
Replace-Text "708÷3=236, 0" "492÷9=54, 6"
Replace-Text "722÷5=144, 2" "508÷6=84, 4"
Replace-Text "981÷5=196, 1" "855÷7=122, 1"
Replace-Text "950÷2=475, 0" "946÷6=157, 4"
Replace-Text "579÷3=193, 0" "420÷5=84, 0"
Replace-Text "343÷6=57, 1" "120÷2=60, 0"

Write-Host "Replacements complete"
